$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 911.6842
$ws.Range("I40").Value = 849.75
$ws.Range("J40").Value = 1017.8571
$ws.Range("K40").Value = 849.75
$ws.Range("L40").Value = 1017.8571
$ws.Range("M40").Value = -674.75
$ws.Range("N40").Value = -1367.8571

$ws.Range("H80").Value = 6583602
$ws.Range("I80").Value = 6250197.5
$ws.Range("J80").Value = 6993946
$ws.Range("K80").Value = 18750592.5
$ws.Range("L80").Value = 20981838
$ws.Range("M80").Value = -18749594.5
$ws.Range("N80").Value = -20983834

$ws.Range("H83").Value = 6583602
$ws.Range("I83").Value = 6250197.5
$ws.Range("J83").Value = 6993946
$ws.Range("K83").Value = 56251777.5
$ws.Range("L83").Value = 62945514
$ws.Range("M83").Value = -56246785.5
$ws.Range("N83").Value = -62955498

$ws.Range("H86").Value = 4167995.8
$ws.Range("I86").Value = 6668272.5
$ws.Range("J86").Value = 867.55554
$ws.Range("K86").Value = 6668272.5
$ws.Range("L86").Value = 867.55554
$ws.Range("M86").Value = -6667149.5
$ws.Range("N86").Value = -3113.55554

$ws.Range("H88").Value = 649.8333
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 649.8333
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 649.8333
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -1461.8333

$ws.Range("H89").Value = 4167995.8
$ws.Range("I89").Value = 6668272.5
$ws.Range("J89").Value = 867.55554
$ws.Range("K89").Value = 33341362.5
$ws.Range("L89").Value = 4337.7777
$ws.Range("M89").Value = -33335746.5
$ws.Range("N89").Value = -15569.7777

$ws.Range("H91").Value = 649.8333
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 649.8333
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 649.8333
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -3457.8333

$ws.Range("H129").Value = 924.9
$ws.Range("I129").Value = 898.4286
$ws.Range("J129").Value = 986.6667
$ws.Range("K129").Value = 2695.2858
$ws.Range("L129").Value = 2960.0001
$ws.Range("M129").Value = 2304.7142
$ws.Range("N129").Value = -12960.0001

$ws.Range("H132").Value = 151883.27
$ws.Range("I132").Value = 2870.691
$ws.Range("J132").Value = 834857.5600000001
$ws.Range("K132").Value = 8612.073
$ws.Range("L132").Value = 2504572.68
$ws.Range("M132").Value = -6082.073
$ws.Range("N132").Value = -2509632.68

$ws.Range("H141").Value = 45363.477
$ws.Range("J141").Value = 77236.086
$ws.Range("L141").Value = 231708.258
$ws.Range("N141").Value = -242068.258

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5929.952
$ws.Range("I32").Value = 4369.763
$ws.Range("K32").Value = 4369.763
$ws.Range("M32").Value = -4082.763

$ws.Range("H61").Value = 1301.5186
$ws.Range("I61").Value = 1201.5769
$ws.Range("J61").Value = 3900
$ws.Range("K61").Value = 1201.5769
$ws.Range("L61").Value = 3900
$ws.Range("M61").Value = -989.5769
$ws.Range("N61").Value = -4324

$ws.Range("H136").Value = 1301.5186
$ws.Range("I136").Value = 1201.5769
$ws.Range("J136").Value = 3900
$ws.Range("K136").Value = 3604.7307
$ws.Range("L136").Value = 11700
$ws.Range("M136").Value = -1054.7307
$ws.Range("N136").Value = -16800

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 260731.48
$ws.Range("I86").Value = 1479.4546
$ws.Range("J86").Value = 1401440.4
$ws.Range("K86").Value = 1479.4546
$ws.Range("L86").Value = 1401440.4
$ws.Range("M86").Value = -356.4546
$ws.Range("N86").Value = -1403686.4

$ws.Range("H89").Value = 260731.48
$ws.Range("I89").Value = 1479.4546
$ws.Range("J89").Value = 1401440.4
$ws.Range("K89").Value = 7397.273
$ws.Range("L89").Value = 7007202
$ws.Range("M89").Value = -1781.273
$ws.Range("N89").Value = -7018434

$ws.Range("H134").Value = 24416386
$ws.Range("I134").Value = 1446.5883
$ws.Range("J134").Value = 143003230
$ws.Range("K134").Value = 4339.7649
$ws.Range("L134").Value = 429009690
$ws.Range("M134").Value = -1804.7649
$ws.Range("N134").Value = -429014760

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2538.9807
$ws.Range("I31").Value = 1365.7576
$ws.Range("K31").Value = 1365.7576
$ws.Range("M31").Value = -1070.7576

$ws.Range("H34").Value = 2538.9807
$ws.Range("I34").Value = 1365.7576
$ws.Range("K34").Value = 1365.7576
$ws.Range("M34").Value = -1163.7576

$ws.Range("H58").Value = 1350.6383
$ws.Range("I58").Value = 1088
$ws.Range("J58").Value = 1705.2
$ws.Range("K58").Value = 1088
$ws.Range("L58").Value = 1705.2
$ws.Range("M58").Value = -885
$ws.Range("N58").Value = -2111.2

$ws.Range("H132").Value = 25267.904
$ws.Range("I132").Value = 38326.37
$ws.Range("J132").Value = 1762.6666
$ws.Range("K132").Value = 114979.11
$ws.Range("L132").Value = 5287.9998
$ws.Range("M132").Value = -112449.11
$ws.Range("N132").Value = -10347.9998

$ws.Range("H134").Value = 1694.9474
$ws.Range("I134").Value = 1564
$ws.Range("J134").Value = 1875
$ws.Range("K134").Value = 4692
$ws.Range("L134").Value = 5625
$ws.Range("M134").Value = -2157
$ws.Range("N134").Value = -10695

$ws.Range("H136").Value = 1350.6383
$ws.Range("I136").Value = 1088
$ws.Range("J136").Value = 1705.2
$ws.Range("K136").Value = 3264
$ws.Range("L136").Value = 5115.6
$ws.Range("M136").Value = -714
$ws.Range("N136").Value = -10215.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1009.3333
$ws.Range("I5").Value = 850
$ws.Range("J5").Value = 1020.7143
$ws.Range("K5").Value = 2550
$ws.Range("L5").Value = 3062.1429
$ws.Range("M5").Value = -2438
$ws.Range("N5").Value = -3286.1429

$ws.Range("H131").Value = 49167336
$ws.Range("J131").Value = 57097496
$ws.Range("L131").Value = 171292488
$ws.Range("N131").Value = -171302568

$ws.Range("H132").Value = 1097.56
$ws.Range("I132").Value = 631.8182
$ws.Range("J132").Value = 1463.5
$ws.Range("K132").Value = 5686.3638
$ws.Range("L132").Value = 13171.5
$ws.Range("M132").Value = -3156.3638
$ws.Range("N132").Value = -18231.5

$ws.Range("H135").Value = 1009.3333
$ws.Range("I135").Value = 850
$ws.Range("J135").Value = 1020.7143
$ws.Range("K135").Value = 7650
$ws.Range("L135").Value = 9186.4287
$ws.Range("M135").Value = -5115
$ws.Range("N135").Value = -14256.4287

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5267379
$ws.Range("I70").Value = 6670720
$ws.Range("J70").Value = 4850
$ws.Range("K70").Value = 6670720
$ws.Range("L70").Value = 4850
$ws.Range("M70").Value = -6670450
$ws.Range("N70").Value = -5390

$ws.Range("H73").Value = 5267379
$ws.Range("I73").Value = 6670720
$ws.Range("J73").Value = 4850
$ws.Range("K73").Value = 6670720
$ws.Range("L73").Value = 4850
$ws.Range("M73").Value = -6669784
$ws.Range("N73").Value = -6722

$ws.Range("H80").Value = 6936.3184
$ws.Range("I80").Value = 2662.375
$ws.Range("J80").Value = 9378.571
$ws.Range("K80").Value = 2662.375
$ws.Range("L80").Value = 9378.571
$ws.Range("M80").Value = -1664.375
$ws.Range("N80").Value = -11374.571

$ws.Range("H83").Value = 6936.3184
$ws.Range("I83").Value = 2662.375
$ws.Range("J83").Value = 9378.571
$ws.Range("K83").Value = 13311.875
$ws.Range("L83").Value = 46892.855
$ws.Range("M83").Value = -8319.875
$ws.Range("N83").Value = -56876.855

$ws.Range("H132").Value = 40851.69
$ws.Range("I132").Value = 1847.5385
$ws.Range("K132").Value = 5542.6155
$ws.Range("M132").Value = -3012.6155

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3418.8948
$ws.Range("I132").Value = 497.06897
$ws.Range("J132").Value = 12833.667
$ws.Range("K132").Value = 1491.20691
$ws.Range("L132").Value = 38501.001
$ws.Range("M132").Value = 1038.79309
$ws.Range("N132").Value = -43561.001
